# Building-Length.docx: replace the embedded chart picture with a plain
# hyperlink run pointing at the image's URL on ura.gov.sg, and leave the
# rest of the document (including the "link" hyperlink further down)
# untouched content-wise.

$d = $word.ActiveDocument

$imageUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/BL1_Framework_Building_Height_Building_Length.png?h=80%25&w=80%25"

# The chart is the document's only inline picture, sitting alone in its
# own BodyText paragraph right after the intro paragraph.
$pic = $d.InlineShapes.Item(1)
$picRange = $pic.Range
$pic.Delete()

# Turn that now-empty paragraph into a hyperlink whose display text is the
# image URL itself (matches how the other external link in this doc is
# built: a Hyperlink-styled run inside a w:hyperlink).
$d.Hyperlinks.Add($picRange, $imageUrl, "", "", $imageUrl)
